$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the trailing run of the first user story ("... hochladen
# koennen, um die Artikel mit anderen Benutzern zu teilen.") so that the
# tail is rebuilt out of several small runs, matching several small manual
# edits made in Word.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
$targetIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $pp = $paras.Item($i)
    if ($pp.Range.Text -like "*hochladen*" -and $pp.Range.Text -like "*anderen Benutzern zu teilen*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $p = $paras.Item($targetIndex)
    $pStart = $p.Range.Start
    $pText = $p.Range.Text

    $oldTail = ", um die Artikel mit anderen Benutzern zu teilen."
    $tailOffset = $pText.IndexOf($oldTail)

    if ($tailOffset -ge 0) {
        $absStart = $pStart + $tailOffset
        $absEnd = $absStart + $oldTail.Length

        # Wipe the old run's text first.
        $wipe = $d.Range($absStart, $absEnd)
        $wipe.Text = ""

        # Re-insert the replacement text as a sequence of distinct runs by
        # inserting each chunk at a fresh collapsed range -- this keeps each
        # chunk as its own <w:r>, mirroring the incremental edits in the
        # authored document.
        $parts = @(
            ", um die",
            "se an",
            " Artikel",
            " anzuhängen und",
            " ",
            "so",
            "mit ",
            " mit ",
            "anderen Benutzern zu teilen."
        )

        $pos = $absStart
        foreach ($part in $parts) {
            $ip = $d.Range($pos, $pos)
            $ip.InsertAfter($part)
            $pos = $pos + $part.Length
        }
    }
}

# ---------------------------------------------------------------------------
# Change 2: add a new list item right after the "Informationsredundanz"
# user story, before the following blank paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
$srcIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $pp = $paras.Item($i)
    if ($pp.Range.Text -like "*Informationsredundanz*") {
        $srcIndex = $i
        break
    }
}

if ($srcIndex -ne -1) {
    $src = $paras.Item($srcIndex)
    $src.Range.InsertParagraphAfter()

    # Re-fetch the freshly created (still empty) paragraph and overwrite it
    # with the exact OOXML we need -- this keeps the paragraph properties
    # minimal (no contextualSpacing override), matching the target markup.
    $paras = $d.Paragraphs
    $newPara = $paras.Item($srcIndex + 1)
    $newRange = $newPara.Range

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr>' +
           '<w:pStyle w:val="Listenabsatz"/>' +
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>' +
           '<w:spacing w:line="288" w:lineRule="auto"/>' +
           '</w:pPr>' +
           '<w:r><w:t>Als Benutzer möchte ich an Artikel angehängte Dokumente herunterladen können, um diese lokal öffnen zu können.</w:t></w:r>' +
           '</w:p>'
    $newRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Change 3 & 4: move <w:lastRenderedPageBreak/> from the run holding "12"
# (estimation table, first column) to the run holding "10".
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$rows = $t.Rows.Count

$row10 = -1
$row12 = -1
for ($r = 1; $r -le $rows; $r++) {
    $c1 = $t.Cell($r, 1)
    $txt = $c1.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "10") { $row10 = $r }
    if ($txt -eq "12") { $row12 = $r }
}

if ($row12 -ne -1) {
    $c12 = $t.Cell($row12, 1)
    $full12 = $c12.Range
    $target12 = $d.Range($full12.Start, $full12.End - 1)
    $target12.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:spacing w:line="288" w:lineRule="auto"/><w:ind w:left="0"/><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:t>12</w:t></w:r></w:p>')
}

if ($row10 -ne -1) {
    $t = $d.Tables.Item(1)
    $c10 = $t.Cell($row10, 1)
    $full10 = $c10.Range
    $target10 = $d.Range($full10.Start, $full10.End - 1)
    $target10.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:spacing w:line="288" w:lineRule="auto"/><w:ind w:left="0"/><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>10</w:t></w:r></w:p>')
}

Write-Host "Done"
